$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step A: append the new roster rows (order matters: it controls the
#     order new strings are interned into the shared-string table) ---
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "Sajvel"

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "J.Mahesh"

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Prahalad"

# "Saif " (row 5) becomes "Md.Saif "
$ws.Cells.Item(5, 2).Value = "Md.Saif "

$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Kamal"

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Akanksha"

# --- Step B: carry the highlight fill from the current Ambuj/Arnab cells
#     (B2:B3) onto the rows they will occupy after reshuffling (B4:B5) ---
$ws.Range("B2:B3").Copy()
$ws.Range("B4:B5").PasteSpecial(-4122)
$ws.Range("B2:B3").ClearFormats()

# --- Step C: reshuffle the names into their final row order ---
$ws.Cells.Item(2, 2).Value = "Md.Saif "
$ws.Cells.Item(3, 2).Value = "Ishaan"
$ws.Cells.Item(4, 2).Value = "Ambuj"
$ws.Cells.Item(5, 2).Value = "Arnab"
$ws.Cells.Item(6, 2).Value = "Nishant"

$ws.Range("E7").Select()
